# Regenerate the "K" column (col G) values for rows 2..74 of Sheet1.
# This corresponds to recalculated/regenerated save_data values
# (commit: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K"), one per row starting at row 2 through row 74.
$kValues = @(
    3, 1, 2, 1, 2, 2, 0, 2, 0, 1,
    1, 2, 0, 1, 2, 1, 0, 2, 0, 3,
    1, 3, 3, 1, 2, 3, 1, 2, 3, 1,
    0, 0, 1, 2, 1, 1, 2, 3, 3, 3,
    1, 1, 3, 3, 3, 4, 2, 0, 1, 1,
    0, 1, 1, 0, 0, 1, 0, 2, 0, 3,
    0, 1, 3, 1, 1, 1, 3, 3, 2, 1,
    0, 1, 1
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}

$endRow = $startRow + $kValues.Length - 1
Write-Host "Updated column G (K) for rows $startRow to $endRow"
